$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "BSHQ240941 / 梁海东" entry (row 16) is dropped from the
# "下周工作安排" table; everything below it shifts up by one row.
$ws.Rows.Item(16).Delete()

# The leftover date-column style (numFmtId 14) used by the populated
# rows' date cells is retired in favour of the format already used for
# the blank placeholder rows (numFmtId 58) -- copy that format across
# the date columns of every populated data row so they share one style.
$ws.Range("E7:H7").Copy()
$ws.Range("E3:H6").PasteSpecial(-4122)
$ws.Range("E15:H16").PasteSpecial(-4122)

# Row 17 used to be a blank placeholder row; it now carries a new
# "下周工作安排" entry, so first give it the same visual format as the
# other populated rows ...
$ws.Range("A16:I16").Copy()
$ws.Range("A17:I17").PasteSpecial(-4122)

# ... then fill in the new entry's text/date content.
$ws.Range("A17").Value = "BSXG250226"
$ws.Range("B17").Value = "公司内部"
$ws.Range("C17").Value = "思路设计"
$ws.Range("D17").Value = "小鼠结肠炎的单细胞测序数据"
$ws.Range("E17").Value = 45715
$ws.Range("F17").Value = 45726

# G17/H17/I17 stay blank, like the un-filled-in "预计完成" / "实际完成" /
# "备注" cells on the other new-entry rows.

# This new entry's description text is longer, so the row is taller
# than its neighbours.
$ws.Rows.Item(17).RowHeight = 42.75
